$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (before) values for the columns that change,
# for rows 2, 3 and 4.
$cols = @("D", "J", "K", "L", "M", "O", "P")

$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value2
    $row3[$col] = $ws.Range($col + "3").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
}

# New values (weekly cycle): row2 <- old row3, row3 <- old row4, row4 <- old row2
foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $row3[$col]
    $ws.Range($col + "3").Value = $row4[$col]
    $ws.Range($col + "4").Value = $row2[$col]
}
